$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT, never letting Excel auto-convert
# numeric-looking strings (e.g. '1.00', '38.80') into real numbers,
# which would silently drop significant trailing/leading zeros.
# Plain-number-looking strings are routed through a scratch cell as a
# text formula result ("=""1.00""") and copied over with values-only
# paste, so no NumberFormat/style churn is introduced on the real cell.
function Set-TextValue($addr, $val) {
    if ($val -match '^\s*-?[0-9]+(\.[0-9]+)?\s*$') {
        $ws.Range("ZZ1").Formula = '="' + $val + '"'
        $ws.Range("ZZ1").Copy()
        $ws.Range($addr).PasteSpecial(-4163)
        $ws.Range("ZZ1").ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

Set-TextValue "D2" '60.704.23'
Set-TextValue "E2" '  +0.45%  '
Set-TextValue "D3" '2.624.97'
Set-TextValue "E3" '  +1.10%  '
Set-TextValue "E4" '  -0.03%  '
Set-TextValue "D5" '585.90'
Set-TextValue "E5" '  +3.62%  '
Set-TextValue "D6" '145.18'
Set-TextValue "E6" '  +1.75%  '
Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  +0.20%  '
Set-TextValue "D8" '0.600'
Set-TextValue "E8" '  +0.18%  '
Set-TextValue "D9" '6.52'
Set-TextValue "E9" '  -0.58%  '
Set-TextValue "E10" '  +0.46%  '
Set-TextValue "E11" '  +1.51%  '
Set-TextValue "D12" '0.154'
Set-TextValue "E12" '  +1.37%  '
Set-TextValue "D13" '3.092.90'
Set-TextValue "E13" '  +1.19%  '
Set-TextValue "D14" '26.12'
Set-TextValue "E14" '  +12.02%  '
Set-TextValue "D15" '60.567.64'
Set-TextValue "E15" '  +0.17%  '
Set-TextValue "E16" '  +0.91%  '
Set-TextValue "D17" '2.631.17'
Set-TextValue "D18" '11.57'
Set-TextValue "E18" '  +3.17%  '
Set-TextValue "D19" '4.73'
Set-TextValue "E19" '  +1.72%  '
Set-TextValue "D20" '349.54'
Set-TextValue "E20" '  +0.90%  '
Set-TextValue "D21" '6.89'
Set-TextValue "E21" '  -1.18%  '
Set-TextValue "D22" '0.998'
Set-TextValue "E22" '  -0.21%  '
Set-TextValue "D23" '0.529'
Set-TextValue "E23" '  -0.75%  '
Set-TextValue "D24" '63.72'
Set-TextValue "E24" '  +0.57%  '
Set-TextValue "D25" '0.998'
Set-TextValue "E25" '  +0.14%  '
Set-TextValue "E26" '  +1.32%  '
Set-TextValue "D27" '8.20'
Set-TextValue "E27" '  +7.14%  '
Set-TextValue "D28" '2.07'
Set-TextValue "E28" '  +15.49%  '
Set-TextValue "D29" '0.0₃0802'
Set-TextValue "E29" '  +2.30%  '
Set-TextValue "D30" '6.57'
Set-TextValue "E30" '  +3.99%  '
Set-TextValue "D31" '170.24'
Set-TextValue "E31" '  +5.74%  '
Set-TextValue "D32" '0.999'
Set-TextValue "E32" '  +0.12%  '
Set-TextValue "D33" '19.57'
Set-TextValue "E33" '  +0.61%  '
Set-TextValue "D34" '4.40'
Set-TextValue "E34" '  +4.10%  '
Set-TextValue "D35" '1.03'
Set-TextValue "E35" '  +6.90%  '
Set-TextValue "D36" '1.32'
Set-TextValue "E36" '  +8.16%  '
Set-TextValue "D37" '1.65'
Set-TextValue "E37" '  +2.78%  '
Set-TextValue "D38" '331.98'
Set-TextValue "E38" '  +12.93%  '
Set-TextValue "B39" 'OKB'
Set-TextValue "C39" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D39" '38.80'
Set-TextValue "E39" '  +2.87%  '
Set-TextValue "B40" 'Filecoin'
Set-TextValue "C40" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D40" '4.00'
Set-TextValue "E40" '  +4.80%  '
Set-TextValue "D41" '0.869'
Set-TextValue "E41" '  +1.24%  '
Set-TextValue "B42" 'Aave'
Set-TextValue "C42" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D42" '134.17'
Set-TextValue "E42" '  -2.22%  '
Set-TextValue "B43" 'RenderToken'
Set-TextValue "C43" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D43" '5.12'
Set-TextValue "E43" '  +5.68%  '
Set-TextValue "D44" '20.19'
Set-TextValue "E44" '  +3.75%  '
Set-TextValue "D45" '0.0998'
Set-TextValue "E45" '  +2.12%  '
Set-TextValue "D46" '1.00'
Set-TextValue "E46" '  +0.35%  '
Set-TextValue "B47" 'Hedera'
Set-TextValue "C47" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D47" '0.0558'
Set-TextValue "E47" '  +2.16%  '
Set-TextValue "B48" 'Mantle'
Set-TextValue "C48" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D48" '0.611'
Set-TextValue "E48" '  +0.95%  '
Set-TextValue "D49" '20.33'
Set-TextValue "E49" '  +3.37%  '
Set-TextValue "E50" '  +2.12%  '
Set-TextValue "D51" '10.73'
Set-TextValue "E51" '  +0.40%  '

$excel.CutCopyMode = $false
